$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# --- "Test Steps" sheet edits -------------------------------------------
# Insert a brand-new column at position 5 (keeps column D's existing width
# untouched) then relocate the old "Action_Keyword" column D data into the
# new column E, leaving column D free for the new "Page Objects" column.
$ws.Columns.Item(5).Insert()
for ($r = 1; $r -le 7; $r++) {
    $v = $ws.Cells.Item($r, 4).Value2
    if ($v -ne $null) {
        $ws.Cells.Item($r, 5).Value = $v
    }
}
$ws.Range("D1:D7").ClearContents()
$ws.Columns.Item(5).ColumnWidth = 16.1796875

# New "Page Objects" header + values for the existing rows
$ws.Range("D1").Value = "Page Objects"
$ws.Range("D4").Value = "txtbx_UserName"
$ws.Range("D5").Value = "txtbx_Password"
$ws.Range("D6").Value = "btn_LogIn"

# Insert a new row above the old "Quit Browser" row for the new
# "Click LogOut Button" test step, and renumber the TS ID of the
# row that gets pushed down.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Login_01"
$ws.Range("B7").Value = "TS_006"

$ws.Range("B8").Value = "TS_007"

$ws.Range("C7").Value = "Click LogOut Button"
$ws.Range("D7").Value = "btn_Signout"
$ws.Range("E7").Value = "doLogout"

[void]$ws.Range("D7").Select()

# --- New "Test Cases" sheet ---------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Test Cases"

$newSheet.Range("A1").Value = "Test Case ID"
$newSheet.Range("B1").Value = "Description"
$newSheet.Range("C1").Value = "Runmode"
$newSheet.Range("A2").Value = "Login_01"
$newSheet.Range("A3").Value = "Login_02"
$newSheet.Range("B2").Value = "Login to the online app"
$newSheet.Range("C2").Value = "Yes"
$newSheet.Range("B3").Value = "Login to the online app"
$newSheet.Range("C3").Value = "No"

$newSheet.Columns.Item(1).ColumnWidth = 12.26953125
$newSheet.Columns.Item(2).ColumnWidth = 22.36328125

[void]$newSheet.Range("D8").Select()

# Keep "Test Steps" as the active sheet/tab (matches the unchanged bookViews
# in the target workbook).
[void]$ws.Activate()
